# Log_of_all_Blogs.xlsx - "Log file updated, with links of Post52"
# Appends a new row (S.No 52) to the Table2 listobject on Sheet1 for the
# "Multi-Instance Resource Allocation Graph | Operating System - M04 P03" post.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$lo = $ws.ListObjects.Item("Table2")

# Grow the table by one row - this expands the table ref / autofilter range
# and the worksheet dimension to B10:F62.
$newRow = $lo.ListRows.Add()

# Carry the formatting (number format / hyperlink look / borders) from the
# last existing data row down onto the freshly added one before filling in
# the values.
$ws.Range("B61:F61").Copy()
$ws.Range("B62:F62").PasteSpecial(-4122)

# Fill in the new row's data. Dev.to link (F) then Title (C) then Hashnode
# link (E) are written in this order so new shared-string entries land in
# the same order as the source edit (158=Dev.to link, 159=Title,
# 160=Hashnode link).
$ws.Range("B62").Value = 52
$ws.Range("F62").Value = "https://dev.to/rahulmishra05/multi-instance-resource-allocation-graph-operating-system-m04-p03-15nh"
$ws.Range("C62").Value = "Multi-Instance Resource Allocation Graph | Operating System - M04 P03"
$ws.Range("D62").Value = 44174
$ws.Range("E62").Value = "https://programmingport.hashnode.dev/multi-instance-resource-allocation-graph-or-operating-system-m04-p03"

# Match the source workbook's updated viewport/selection (scrolled so column
# D is at the left edge, cell E62 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 4
$null = $ws.Range("E62").Select()
